$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($addr, $val) {
    $ws.Range($addr).Value = $val
}

# Row 2 - Bitcoin
Set-Cell "D2" "68.612.04"
Set-Cell "E2" "  +1.07%  "

# Row 3 - Ethereum
Set-Cell "D3" "3.867.18"
Set-Cell "E3" "  +0.29%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  -0.06%  "

# Row 5 - BNB
Set-Cell "D5" "603.25"
Set-Cell "E5" "  +0.82%  "

# Row 6 - Solana
Set-Cell "D6" "173.03"
Set-Cell "E6" "  +4.06%  "

# Row 7 - LidoStakedEther
Set-Cell "D7" "3.865.58"
Set-Cell "E7" "  +0.31%  "

# Row 8 - USDC
Set-Cell "E8" "  -0.06%  "

# Row 9 - XRP
Set-Cell "E9" "  +1.18%  "

# Row 10 - Dogecoin
Set-Cell "D10" "0.170"
Set-Cell "E10" "  +3.32%  "

# Row 11 - Toncoin
Set-Cell "D11" "6.54"
Set-Cell "E11" "  +3.66%  "

# Row 12/13 - Cardano <-> ShibaInu swap
Set-Cell "B12" "ShibaInu"
Set-Cell "C12" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-Cell "D12" "0.0000293"
Set-Cell "E12" "  +18.27%  "

Set-Cell "B13" "Cardano"
Set-Cell "C13" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-Cell "D13" "0.464"
Set-Cell "E13" "  +1.76%  "

# Row 14 - Avalanche
Set-Cell "D14" "37.38"
Set-Cell "E14" "  +1.55%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Cell "D15" "4.515.36"
Set-Cell "E15" "  +0.25%  "

# Row 16 - WrappedEther
Set-Cell "D16" "3.877.87"
Set-Cell "E16" "  +0.75%  "

# Row 17 - WrappedBTC
Set-Cell "D17" "68.674.12"
Set-Cell "E17" "  +1.14%  "

# Row 18/19 - Polkadot <-> Chainlink swap
Set-Cell "B18" "Chainlink"
Set-Cell "C18" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-Cell "D18" "18.35"
Set-Cell "E18" "  +1.08%  "

Set-Cell "B19" "Polkadot"
Set-Cell "C19" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-Cell "D19" "7.51"
Set-Cell "E19" "  +1.62%  "

# Row 20 - TRON
Set-Cell "D20" "0.111"
Set-Cell "E20" "  +0.82%  "

# Row 21 - Uniswap
Set-Cell "E21" "  +2.25%  "

# Row 22 - BitcoinCash
Set-Cell "D22" "473.07"
Set-Cell "E22" "  +1.56%  "

# Row 23 - Polygon
Set-Cell "D23" "0.736"
Set-Cell "E23" "  +1.09%  "

# Row 24 - PEPE
Set-Cell "E24" "  +1.58%  "

# Row 25 - Litecoin
Set-Cell "E25" "  +0.91%  "

# Row 26 - Fetch.AI
Set-Cell "D26" "2.30"
Set-Cell "E26" "  +3.58%  "

# Row 27 - InternetComputer(DFINITY)
Set-Cell "D27" "12.33"
Set-Cell "E27" "  +1.83%  "

# Row 28 - RenderToken
Set-Cell "D28" "10.57"
Set-Cell "E28" "  +5.90%  "

# Row 29 - Dai
Set-Cell "E29" "  +0.04%  "

# Row 30 - PancakeSwap
Set-Cell "E30" "  +0.45%  "

# Row 31 - WrappedeETH
Set-Cell "D31" "4.017.81"
Set-Cell "E31" "  +0.25%  "

# Row 32 - NEARProtocol
Set-Cell "D32" "7.82"
Set-Cell "E32" "  +1.89%  "

# Row 33 - ImmutableX
Set-Cell "E33" "  +1.28%  "

# Row 34 - EthereumClassic
Set-Cell "E34" "  +1.26%  "

# Row 35 - Aptos
Set-Cell "D35" "9.50"
Set-Cell "E35" "  +1.16%  "

# Row 36 - RenzoRestakedETH
Set-Cell "D36" "3.833.14"
Set-Cell "E36" "  -0.01%  "

# Row 37 - dogwifhat
Set-Cell "D37" "3.95"
Set-Cell "E37" "  +19.77%  "

# Row 38 - Hedera
Set-Cell "E38" "  +2.04%  "

# Row 39 - Filecoin
Set-Cell "E39" "  +2.68%  "

# Row 40/41 - Kaspa <-> Mantle swap
Set-Cell "B40" "Mantle"
Set-Cell "C40" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-Cell "D40" "1.02"
Set-Cell "E40" "  +0.99%  "

Set-Cell "B41" "Kaspa"
Set-Cell "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-Cell "D41" "0.140"
Set-Cell "E41" "  +0.10%  "

# Row 42 - FirstDigitalUSD
Set-Cell "E42" "  -0.01%  "

# Row 43 - TheGraph
Set-Cell "D43" "0.323"
Set-Cell "E43" "  +3.70%  "

# Row 44 - FLOKI
Set-Cell "D44" "0.000302"
Set-Cell "E44" "  +11.28%  "

# Row 45 - Stacks
Set-Cell "E45" "  +1.54%  "

# Row 46/47/48 - Cosmos, USDe, Bittensor 3-way rotate
Set-Cell "B46" "Bittensor"
Set-Cell "C46" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-Cell "D46" "424.08"
Set-Cell "E46" "  -0.79%  "

Set-Cell "B47" "Cosmos"
Set-Cell "C47" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-Cell "D47" "8.81"
Set-Cell "E47" "  +3.40%  "

Set-Cell "B48" "USDe"
Set-Cell "C48" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-Cell "D48" "1.00"
Set-Cell "E48" "  -0.02%  "

# Row 49 - OKB
Set-Cell "D49" "46.69"
Set-Cell "E49" "  -1.10%  "

# Row 50 - VeChain
Set-Cell "D50" "0.0363"
Set-Cell "E50" "  +2.88%  "

# Row 51 - Monero
Set-Cell "D51" "142.26"
Set-Cell "E51" "  -0.68%  "
